$d = $word.ActiveDocument

$pairs = @(
    @("63×22=", "88×63="),
    @("71×37=", "31×82="),
    @("91×15=", "92×73="),
    @("17×86=", "33×61="),
    @("73×75=", "65×78="),
    @("72×68=", "46×11="),
    @("85×33=", "17×19="),
    @("45×17=", "53×99="),
    @("80×13=", "36×12="),
    @("78×97=", "42×80="),
    @("79×84=", "34×48="),
    @("62×69=", "32×84="),
    @("59×26=", "59×73="),
    @("44×27=", "30×61="),
    @("86×45=", "20×82="),
    @("66×14=", "88×74="),
    @("62×84=", "90×54="),
    @("64×31=", "50×21="),
    @("12×72=", "54×52="),
    @("31×33=", "59×25="),
    @("85×83=", "79×76="),
    @("57×28=", "97×14="),
    @("53×88=", "88×86="),
    @("95×13=", "96×57="),
    @("87×60=", "91×46=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
